$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 114.30769
$ws.Range("I28").Value = 114.30769
$ws.Range("K28").Value = 114.30769
$ws.Range("M28").Value = 370.69231
$ws.Range("H62").Value = 2295.4443
$ws.Range("I62").Value = 2293.3333
$ws.Range("J62").Value = 2299.6667
$ws.Range("K62").Value = 2293.3333
$ws.Range("L62").Value = 2299.6667
$ws.Range("M62").Value = -1669.3333
$ws.Range("N62").Value = -3547.6667
$ws.Range("H65").Value = 2295.4443
$ws.Range("I65").Value = 2293.3333
$ws.Range("J65").Value = 2299.6667
$ws.Range("K65").Value = 11466.6665
$ws.Range("L65").Value = 11498.3335
$ws.Range("M65").Value = -8346.666499999999
$ws.Range("N65").Value = -17738.3335
$ws.Range("H111").Value = 2176.6667
$ws.Range("I111").Value = 1765
$ws.Range("J111").Value = 3000
$ws.Range("K111").Value = 5295
$ws.Range("L111").Value = 9000
$ws.Range("M111").Value = -2228
$ws.Range("N111").Value = -15134
$ws.Range("H113").Value = 2344.9
$ws.Range("I113").Value = 1993.3334
$ws.Range("J113").Value = 3399.6
$ws.Range("K113").Value = 1993.3334
$ws.Range("L113").Value = 3399.6
$ws.Range("M113").Value = 1260.6666
$ws.Range("N113").Value = -9907.6
$ws.Range("H116").Value = 2620
$ws.Range("I116").Value = 1533.3334
$ws.Range("J116").Value = 2891.6667
$ws.Range("K116").Value = 1533.3334
$ws.Range("L116").Value = 2891.6667
$ws.Range("M116").Value = 1908.6666
$ws.Range("N116").Value = -9775.6667
$ws.Range("H135").Value = 1808.2354
$ws.Range("I135").Value = 1833.1666
$ws.Range("J135").Value = 1748.4
$ws.Range("K135").Value = 16498.4994
$ws.Range("L135").Value = 15735.6
$ws.Range("M135").Value = -13963.4994
$ws.Range("N135").Value = -20805.6
$ws.Range("H138").Value = 1720.1052
$ws.Range("I138").Value = 1174.421
$ws.Range("J138").Value = 2265.7896
$ws.Range("K138").Value = 3523.263
$ws.Range("L138").Value = 6797.3688
$ws.Range("M138").Value = 1616.737
$ws.Range("N138").Value = -17077.3688

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1186.4166
$ws.Range("I2").Value = 951
$ws.Range("J2").Value = 1385.6154
$ws.Range("K2").Value = 951
$ws.Range("L2").Value = 1385.6154
$ws.Range("M2").Value = -838
$ws.Range("N2").Value = -1611.6154
$ws.Range("H116").Value = 1186.4166
$ws.Range("I116").Value = 951
$ws.Range("J116").Value = 1385.6154
$ws.Range("K116").Value = 951
$ws.Range("L116").Value = 1385.6154
$ws.Range("M116").Value = 1343
$ws.Range("N116").Value = -5973.6154
$ws.Range("H122").Value = 2421.9443
$ws.Range("I122").Value = 2066.6667
$ws.Range("J122").Value = 3132.5
$ws.Range("K122").Value = 6200.000100000001
$ws.Range("L122").Value = 9397.5
$ws.Range("M122").Value = -3750.000100000001
$ws.Range("N122").Value = -14297.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1186.4166
$ws.Range("I3").Value = 951
$ws.Range("J3").Value = 1385.6154
$ws.Range("K3").Value = 951
$ws.Range("L3").Value = 1385.6154
$ws.Range("M3").Value = -837
$ws.Range("N3").Value = -1613.6154
$ws.Range("H94").Value = 943.4400000000001
$ws.Range("I94").Value = 864.4706
$ws.Range("J94").Value = 1111.25
$ws.Range("K94").Value = 864.4706
$ws.Range("L94").Value = 1111.25
$ws.Range("M94").Value = -413.4706
$ws.Range("N94").Value = -2013.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1373.3334
$ws.Range("I16").Value = 860
$ws.Range("J16").Value = 1630
$ws.Range("K16").Value = 860
$ws.Range("L16").Value = 1630
$ws.Range("M16").Value = -573
$ws.Range("N16").Value = -2204
$ws.Range("H99").Value = 1550.9697
$ws.Range("I99").Value = 993.3333
$ws.Range("J99").Value = 2015.6666
$ws.Range("K99").Value = 993.3333
$ws.Range("L99").Value = 2015.6666
$ws.Range("M99").Value = 504.6667
$ws.Range("N99").Value = -5011.6666
$ws.Range("H107").Value = 616
$ws.Range("I107").Value = 561.75
$ws.Range("J107").Value = 833
$ws.Range("K107").Value = 561.75
$ws.Range("L107").Value = 833
$ws.Range("M107").Value = 1358.25
$ws.Range("N107").Value = -4673
$ws.Range("H113").Value = 1373.3334
$ws.Range("I113").Value = 860
$ws.Range("J113").Value = 1630
$ws.Range("K113").Value = 860
$ws.Range("L113").Value = 1630
$ws.Range("M113").Value = 1310
$ws.Range("N113").Value = -5970
$ws.Range("H126").Value = 1550.9697
$ws.Range("I126").Value = 993.3333
$ws.Range("J126").Value = 2015.6666
$ws.Range("K126").Value = 2979.9999
$ws.Range("L126").Value = 6046.9998
$ws.Range("M126").Value = -509.9998999999998
$ws.Range("N126").Value = -10986.9998

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H68").Value = 386.66666
$ws.Range("I68").Value = 364
$ws.Range("J68").Value = 500
$ws.Range("K68").Value = 1092
$ws.Range("L68").Value = 1500
$ws.Range("M68").Value = -281
$ws.Range("N68").Value = -3122
$ws.Range("H71").Value = 386.66666
$ws.Range("I71").Value = 364
$ws.Range("J71").Value = 500
$ws.Range("K71").Value = 3276
$ws.Range("L71").Value = 4500
$ws.Range("M71").Value = 780
$ws.Range("N71").Value = -12612
$ws.Range("H97").Value = 280.2143
$ws.Range("I97").Value = 269.22223
$ws.Range("K97").Value = 807.66669
$ws.Range("M97").Value = -311.66669

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 9800
$ws.Range("J26").Value = 9800
$ws.Range("L26").Value = 9800
$ws.Range("N26").Value = -10360
$ws.Range("H50").Value = 9800
$ws.Range("J50").Value = 9800
$ws.Range("L50").Value = 9800
$ws.Range("N50").Value = -10796
$ws.Range("H102").Value = 2002.6
$ws.Range("I102").Value = 1999
$ws.Range("J102").Value = 2005
$ws.Range("K102").Value = 1999
$ws.Range("L102").Value = 2005
$ws.Range("M102").Value = -377
$ws.Range("N102").Value = -5249
$ws.Range("H107").Value = 718.3333
$ws.Range("I107").Value = 630.3125
$ws.Range("K107").Value = 630.3125
$ws.Range("M107").Value = 1289.6875
$ws.Range("H113").Value = 8334138.5
$ws.Range("I113").Value = 17857738
$ws.Range("J113").Value = 989.375
$ws.Range("K113").Value = 17857738
$ws.Range("L113").Value = 989.375
$ws.Range("M113").Value = -17855568
$ws.Range("N113").Value = -5329.375
$ws.Range("H122").Value = 2073.6843
$ws.Range("I122").Value = 1440
$ws.Range("J122").Value = 2777.7778
$ws.Range("K122").Value = 4320
$ws.Range("L122").Value = 8333.3334
$ws.Range("M122").Value = -1870
$ws.Range("N122").Value = -13233.3334

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1199.4736
$ws.Range("I61").Value = 1026
$ws.Range("J61").Value = 1850
$ws.Range("K61").Value = 1026
$ws.Range("L61").Value = 1850
$ws.Range("M61").Value = -824
$ws.Range("N61").Value = -2254
$ws.Range("H93").Value = 2386.7144
$ws.Range("I93").Value = 2567.1667
$ws.Range("J93").Value = 1304
$ws.Range("K93").Value = 2567.1667
$ws.Range("L93").Value = 1304
$ws.Range("M93").Value = -1319.1667
$ws.Range("N93").Value = -3800
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H113").Value = 1199.4736
$ws.Range("I113").Value = 1026
$ws.Range("J113").Value = 1850
$ws.Range("K113").Value = 1026
$ws.Range("L113").Value = 1850
$ws.Range("M113").Value = 1144
$ws.Range("N113").Value = -6190
$ws.Range("H122").Value = 2553
$ws.Range("I122").Value = 2496.4614
$ws.Range("K122").Value = 7489.3842
$ws.Range("M122").Value = -5039.3842

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 212.5
$ws.Range("I107").Value = 218.75
$ws.Range("K107").Value = 656.25
$ws.Range("M107").Value = 1263.75
$ws.Range("H113").Value = 519.6
$ws.Range("I113").Value = 542.5833
$ws.Range("K113").Value = 1627.7499
$ws.Range("M113").Value = 542.2501
